$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrement the "剩余" (remaining) value in column E by 1 for every data
# row (2-99), except row 36 which keeps its original value.
for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current - 1
    }
}
